$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-5, columns A and B
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 156

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 143

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 93

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 77

# Remove row 6 entirely (data shrank from A1:B6 to A1:B5)
$ws.Range("A6:B6").Delete()
